$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row at position 5, shifting the existing rows 5-14 down to 6-15
$ws.Rows.Item(5).Insert()

# Fill in the new row with the LP solver setting (label in A, value in B)
$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"

# Match the label formatting used by the other option rows in column A
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Make "general" the active sheet and select the newly added row
$ws.Activate()
$ws.Range("A5:B5").Select()
